$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I = "I0", J = "IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold font, border, centered) used by the
# other header cells (e.g. H1) by copying formats only.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for I2:J34
$iValues = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 5, 1, 1, 4, 2)
$jValues = @(4, 6, 6, 6, 4, 2, 2, 6, 6, 6, 4, 6, 6, 6, 6, 7, 4, 8, 7, 5, 6, 7, 4, 6, 4, 5, 6, 5, 6, 6, 5, 6, 2)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
